# Atualiza os rótulos da primeira linha de cada planilha para que o
# Power BI identifique automaticamente a linha como cabeçalho.

$wb = $excel.ActiveWorkbook

# Planilhas cujo cabeçalho usa o padrão "Ano <ano>"
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Value2
    $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Value2
    $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Value2
    $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Value2
}

# Planilha com intervalos de anos usa o padrão "Intervalo <intervalo>"
$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIntervalo.Range("B1").Value = "Intervalo " + $wsIntervalo.Range("B1").Value2
$wsIntervalo.Range("C1").Value = "Intervalo " + $wsIntervalo.Range("C1").Value2
$wsIntervalo.Range("D1").Value = "Intervalo " + $wsIntervalo.Range("D1").Value2
$wsIntervalo.Range("E1").Value = "Intervalo " + $wsIntervalo.Range("E1").Value2

# Planilha com apenas uma coluna de dado (ano) usa o padrão "Ano <ano>"
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano " + $wsCusto.Range("B1").Value2
